$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values remain text (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.007.56'
$ws.Range("E2").Value = '  -2.35%  '

$ws.Range("D3").Value = '1.888.09'
$ws.Range("E3").Value = '  -2.93%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").Value = '329.77'
$ws.Range("E5").Value = '  -3.39%  '

$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.25%  '

$ws.Range("D7").Value = '0.4576'
$ws.Range("E7").Value = '  -4.20%  '

$ws.Range("D8").Value = '0.4108'
$ws.Range("E8").Value = '  -0.88%  '

$ws.Range("D9").Value = '47.68'
$ws.Range("E9").Value = '  -1.82%  '

$ws.Range("D10").Value = '0.07956'
$ws.Range("E10").Value = '  -3.56%  '

$ws.Range("D11").Value = '0.9959'
$ws.Range("E11").Value = '  -4.39%  '

$ws.Range("D12").Value = '21.73'
$ws.Range("E12").Value = '  -3.97%  '

$ws.Range("D13").Value = '1.901.97'
$ws.Range("E13").Value = '  -2.90%  '

$ws.Range("D14").Value = '5.915'
$ws.Range("E14").Value = '  -4.22%  '

$ws.Range("D15").Value = '7.077'
$ws.Range("E15").Value = '  -4.76%  '

$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.41%  '

$ws.Range("D17").Value = '88.63'
$ws.Range("E17").Value = '  -4.34%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '0.06564'
$ws.Range("E18").Value = '  -1.87%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.00001025'
$ws.Range("E19").Value = '  -3.65%  '

$ws.Range("E20").Value = '  -3.37%  '

$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.05%  '

$ws.Range("D22").Value = '28.990.76'
$ws.Range("E22").Value = '  -2.41%  '

$ws.Range("E23").Value = '  -3.19%  '

$ws.Range("D24").Value = '11.41'
$ws.Range("E24").Value = '  +1.24%  '

$ws.Range("E25").Value = '  -3.72%  '

$ws.Range("D26").Value = '2.124.75'
$ws.Range("E26").Value = '  -2.96%  '

$ws.Range("D27").Value = '156.01'
$ws.Range("E27").Value = '  -3.29%  '

$ws.Range("D28").Value = '19.58'
$ws.Range("E28").Value = '  -2.84%  '

$ws.Range("D29").Value = '2.085'
$ws.Range("E29").Value = '  -4.83%  '

$ws.Range("D30").Value = '5.495'
$ws.Range("E30").Value = '  -2.35%  '

$ws.Range("D31").Value = '117.47'
$ws.Range("E31").Value = '  -4.00%  '

$ws.Range("D32").Value = '1.041'
$ws.Range("E32").Value = '  +1.40%  '

$ws.Range("D33").Value = '0.09319'
$ws.Range("E33").Value = '  -3.39%  '

$ws.Range("D34").Value = '1.407'
$ws.Range("E34").Value = '  -4.55%  '

$ws.Range("D35").Value = '3.522'
$ws.Range("E35").Value = '  -4.42%  '

$ws.Range("D36").Value = '5.292'
$ws.Range("E36").Value = '  -3.59%  '

$ws.Range("D37").Value = '0.06051'
$ws.Range("E37").Value = '  -3.65%  '

$ws.Range("D38").Value = '0.02226'
$ws.Range("E38").Value = '  -4.17%  '

$ws.Range("D39").Value = '8.373'
$ws.Range("E39").Value = '  -3.88%  '

$ws.Range("D40").Value = '1.171'
$ws.Range("E40").Value = '  -1.82%  '

$ws.Range("D41").Value = '0.9991'
$ws.Range("E41").Value = '  -0.14%  '

$ws.Range("D42").Value = '0.5783'
$ws.Range("E42").Value = '  -5.15%  '

$ws.Range("D43").Value = '0.1822'
$ws.Range("E43").Value = '  -4.38%  '

$ws.Range("E44").Value = '  -5.39%  '

$ws.Range("D45").Value = '1.259'
$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("D46").Value = '0.07497'
$ws.Range("E46").Value = '  +0.83%  '

$ws.Range("D47").Value = '2.283'
$ws.Range("E47").Value = '  -1.92%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '0.5455'
$ws.Range("E48").Value = '  -4.43%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '11.94'
$ws.Range("E49").Value = '  -5.28%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.901'
$ws.Range("E50").Value = '  -4.55%  '

$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '111.26'
$ws.Range("E51").Value = '  -3.16%  '
